# Traditional -> Simplified Chinese translation update, plus a couple of
# English re-translations, for "Email 5-1 [TEMPLATE] Partner email – invite
# revoked" (zh Crowdin translation).

$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, $wdFindContinue, $false, $replace, `
                             $wdReplaceAll) | Out-Null
}

# Top banner language links
Replace-Text "英語" "英语"
Replace-Text " / 葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語" " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语"

# Table: "Brief" heading + body copy
Replace-Text "簡介" "简要"
Replace-Text "發送給在目標國家的合作夥伴的電子郵件，這些合作夥伴已回應參加，但在截止日期前未提交文件。 我們將取消他們的邀請。 將通過 customer.io 發送" "发给在目标国家已确认出席但未在截止日期前提交文件的合作伙伴的邮件。 我们将撤回他们的邀请。 将通过 customer.io 发送"

# Table: "Target audience" heading + body copy
Replace-Text "目標受眾" "目标受众"
Replace-Text "未按時提交文件的被邀請合作夥伴" "未及时提交文件的被邀请合作伙伴"

# Subject line
Replace-Text "主題行" "主题行"
Replace-Text "[活動名稱]" "[活动名称]"
Replace-Text " 註冊" " 注册"

# Heading
Replace-Text "沒有及時收到您的文件" "没有及时收到您的文件"

# Greeting placeholder
Replace-Text "[合作夥伴姓名]" "[PARTNER NAME]"

# Deadline paragraph
Replace-Text "截止日期（" "We didn’t receive your documents by the deadline ("
Replace-Text "[日月年]" "[DD Mmm YYYY]"
Replace-Text "）前沒有收到您的文件。 很遺憾，無法為您辦理 " "). 很遗憾，无法为您办理 "
Replace-Text " 的註冊手續。" " 的注册手续。"

# Well-wishes paragraph (includes commented span)
Replace-Text "衷心祝愿您一切順利，並希望在下一次 " "衷心祝愿您一切顺利，并希望在下一次 "
Replace-Text "會議/研討會/聯盟會員旅行" "会议/研讨会/联盟会员旅行"
Replace-Text "中見到您。" "中见到您。"

# Contact paragraph
Replace-Text "如有任何疑問，請通過 " "如有任何疑问，请通过 "
Replace-Text "[電子郵件地址]" "[电子邮件地址]"
Replace-Text "[WHATSAPP 號碼]" "[WHATSAPP 号码]"
Replace-Text " (WhatsApp) 聯繫您的區域經理 " " (WhatsApp) 联系您的区域经理 "

# Comments ("Choose one") — Find doesn't reach into the comments story, so
# update the comment text directly via the Comments collection.
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $comment = $d.Comments.Item($i)
    if ($comment.Range.Text -eq "選擇其中一個") {
        $comment.Range.Text = "选择任一"
    }
}
